{"js": "// The diff appends three sibling runs \u2014 \" (\", \"Changed main\", \")\" \u2014 right\n// after the existing run in the document's first paragraph\n// (\"This is a Microsoft word document.\"), so the visible text becomes\n// \"This is a Microsoft word document. (Changed main)\" while keeping the\n// new text in three separate <w:r> elements (not merged into the\n// original run, and not merged with each other).\n//\n// A plain paragraph.insertText(...) call (even invoked three times) gets\n// coalesced by the engine into a single run because the runs would share\n// identical formatting, so instead we splice in literal OOXML (wrapped in\n// the required flat-OPC envelope) via Range.insertOoxml, which preserves\n// run boundaries exactly as authored.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Collapse to a range right at the end of the paragraph's text (i.e.\n// just before the paragraph mark), so the inserted runs land inside the\n// existing paragraph instead of spawning a new one.\nconst insertionPoint = firstParagraph.getRange(Word.RangeLocation.end);\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  \"<w:r><w:t>Changed main</w:t></w:r>\" +\n  \"<w:r><w:t>)</w:t></w:r>\" +\n  \"</w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\ninsertionPoint.insertOoxml(flatOpcXml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# The diff appends three sibling runs - \" (\", \"Changed main\", \")\" - right\n# after the existing run in the document's first paragraph\n# (\"This is a Microsoft word document.\"), so the visible text becomes\n# \"This is a Microsoft word document. (Changed main)\" while keeping the\n# new text in three separate <w:r> elements (not merged into the\n# original run, and not merged with each other).\n#\n# Range.InsertAfter (or setting .Text) would coalesce the new text into\n# the existing run, so instead we splice in literal OOXML (wrapped in the\n# required flat-OPC envelope) via Range.InsertXML, which preserves run\n# boundaries exactly as authored.\n\n$d = $word.ActiveDocument\n$firstParagraph = $d.Paragraphs(1)\n$r = $firstParagraph.Range\n\n# Exclude the trailing paragraph mark, then collapse to its end so the\n# inserted runs land inside the existing paragraph instead of spawning a\n# new one.\n$r.MoveEnd(1, -1)\n$r.Collapse(0)\n\n$flatOpcXml = '<?xml version=\"1.0\" standalone=\"yes\"?>' + `\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n  '<pkg:xmlData>' + `\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n  '<w:body><w:p>' + `\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' + `\n  '<w:r><w:t>Changed main</w:t></w:r>' + `\n  '<w:r><w:t>)</w:t></w:r>' + `\n  '</w:p></w:body></w:document>' + `\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($flatOpcXml, \"End\")\n"}
